$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 2194.6667
$ws.Range("I53").Value = 2784.111
$ws.Range("J53").Value = 426.33334
$ws.Range("K53").Value = 2784.111
$ws.Range("L53").Value = 426.33334
$ws.Range("M53").Value = -2147.111
$ws.Range("N53").Value = -1700.33334
$ws.Range("H113").Value = 12612.7
$ws.Range("I113").Value = 14934
$ws.Range("J113").Value = 3327.5
$ws.Range("K113").Value = 14934
$ws.Range("L113").Value = 3327.5
$ws.Range("M113").Value = -11680
$ws.Range("N113").Value = -9835.5
$ws.Range("H137").Value = 1722.3704
$ws.Range("I137").Value = 1710.7368
$ws.Range("J137").Value = 1750
$ws.Range("K137").Value = 5132.2104
$ws.Range("L137").Value = 5250
$ws.Range("M137").Value = -2582.2104
$ws.Range("N137").Value = -10350
$ws.Range("H138").Value = 7691.5
$ws.Range("I138").Value = 5861.091
$ws.Range("J138").Value = 8182.5854
$ws.Range("K138").Value = 17583.273
$ws.Range("L138").Value = 24547.7562
$ws.Range("M138").Value = -12443.273
$ws.Range("N138").Value = -34827.7562
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 805.6799999999999
$ws.Range("I32").Value = 775.8333
$ws.Range("J32").Value = 1074.3
$ws.Range("K32").Value = 775.8333
$ws.Range("L32").Value = 1074.3
$ws.Range("M32").Value = -488.8333
$ws.Range("N32").Value = -1648.3
$ws.Range("H61").Value = 1566.2307
$ws.Range("I61").Value = 1530.0834
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1530.0834
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1318.0834
$ws.Range("N61").Value = -2424
$ws.Range("H136").Value = 1566.2307
$ws.Range("I136").Value = 1530.0834
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 4590.2502
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -2040.2502
$ws.Range("N136").Value = -11100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5775.231
$ws.Range("I134").Value = 6614.5
$ws.Range("J134").Value = 1938.5714
$ws.Range("K134").Value = 19843.5
$ws.Range("L134").Value = 5815.7142
$ws.Range("M134").Value = -17308.5
$ws.Range("N134").Value = -10885.7142
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1687.25
$ws.Range("I22").Value = 1499.6
$ws.Range("K22").Value = 1499.6
$ws.Range("M22").Value = -1149.6
$ws.Range("H31").Value = 1868.9219
$ws.Range("I31").Value = 906.4857
$ws.Range("J31").Value = 3030.4827
$ws.Range("K31").Value = 906.4857
$ws.Range("L31").Value = 3030.4827
$ws.Range("M31").Value = -611.4857
$ws.Range("N31").Value = -3620.4827
$ws.Range("H34").Value = 1868.9219
$ws.Range("I34").Value = 906.4857
$ws.Range("J34").Value = 3030.4827
$ws.Range("K34").Value = 906.4857
$ws.Range("L34").Value = 3030.4827
$ws.Range("M34").Value = -704.4857
$ws.Range("N34").Value = -3434.4827
$ws.Range("H58").Value = 1977186
$ws.Range("I58").Value = 3106508.5
$ws.Range("K58").Value = 3106508.5
$ws.Range("M58").Value = -3106305.5
$ws.Range("H132").Value = 1503.7959
$ws.Range("I132").Value = 871.19446
$ws.Range("J132").Value = 3255.6155
$ws.Range("K132").Value = 2613.58338
$ws.Range("L132").Value = 9766.8465
$ws.Range("M132").Value = -83.58338000000003
$ws.Range("N132").Value = -14826.8465
$ws.Range("H134").Value = 1397.9487
$ws.Range("I134").Value = 1141.1562
$ws.Range("K134").Value = 3423.4686
$ws.Range("M134").Value = -888.4685999999997
$ws.Range("H136").Value = 1977186
$ws.Range("I136").Value = 3106508.5
$ws.Range("K136").Value = 9319525.5
$ws.Range("M136").Value = -9316975.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3321.4285
$ws.Range("I68").Value = 1425
$ws.Range("J68").Value = 4080
$ws.Range("K68").Value = 4275
$ws.Range("L68").Value = 12240
$ws.Range("M68").Value = -3464
$ws.Range("N68").Value = -13862
$ws.Range("H71").Value = 3321.4285
$ws.Range("I71").Value = 1425
$ws.Range("J71").Value = 4080
$ws.Range("K71").Value = 12825
$ws.Range("L71").Value = 36720
$ws.Range("M71").Value = -8769
$ws.Range("N71").Value = -44832
$ws.Range("H107").Value = 4000
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H122").Value = 1088.9166
$ws.Range("J122").Value = 1269.2858
$ws.Range("L122").Value = 11423.5722
$ws.Range("N122").Value = -16323.5722
$ws.Range("H131").Value = 25036458
$ws.Range("I131").Value = 50000520
$ws.Range("J131").Value = 72396.2
$ws.Range("K131").Value = 150001560
$ws.Range("L131").Value = 217188.6
$ws.Range("M131").Value = -149996520
$ws.Range("N131").Value = -227268.6
$ws.Range("H137").Value = 4098.143
$ws.Range("I137").Value = 1333.625
$ws.Range("J137").Value = 5799.385
$ws.Range("K137").Value = 4000.875
$ws.Range("L137").Value = 17398.155
$ws.Range("M137").Value = 1099.125
$ws.Range("N137").Value = -27598.155
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2016.65
$ws.Range("I102").Value = 2019.9231
$ws.Range("J102").Value = 2010.5714
$ws.Range("K102").Value = 2019.9231
$ws.Range("L102").Value = 2010.5714
$ws.Range("M102").Value = -397.9231
$ws.Range("N102").Value = -5254.5714
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2796.9285
$ws.Range("I22").Value = 1096.7142
$ws.Range("J22").Value = 4497.143
$ws.Range("K22").Value = 1096.7142
$ws.Range("L22").Value = 4497.143
$ws.Range("M22").Value = -801.7141999999999
$ws.Range("N22").Value = -5087.143
$ws.Range("H27").Value = 2796.9285
$ws.Range("I27").Value = 1096.7142
$ws.Range("J27").Value = 4497.143
$ws.Range("K27").Value = 1096.7142
$ws.Range("L27").Value = 4497.143
$ws.Range("M27").Value = -989.7141999999999
$ws.Range("N27").Value = -4711.143
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2777.8
$ws.Range("J96").Value = 3224.75
$ws.Range("L96").Value = 3224.75
$ws.Range("N96").Value = -5970.75
$ws.Range("H132").Value = 1466.6522
$ws.Range("I132").Value = 960.94116
$ws.Range("J132").Value = 2899.5
$ws.Range("K132").Value = 2882.82348
$ws.Range("L132").Value = 8698.5
$ws.Range("M132").Value = -352.82348
$ws.Range("N132").Value = -13758.5
$ws.Range("H136").Value = 14247347
$ws.Range("I136").Value = 27780196
$ws.Range("J136").Value = 2243.9473
$ws.Range("K136").Value = 83340588
$ws.Range("L136").Value = 6731.841899999999
$ws.Range("M136").Value = -83338038
